$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.125.82"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5212"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07558"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.680.60"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.421"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5443"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008037"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.48"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.159.10"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.743"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.57"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.244"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.48"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1237"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.481"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.10%  "

$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06313"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.373"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.05%  "

$ws.Range("E30").Value = "  -1.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.507"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("E32").Value = "  -3.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.649"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.003"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6005"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.763"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.116.07"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.41%  "

$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.060"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8640"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.63"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.822.47"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("E46").Value = "  -2.65%  "

$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.058"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05253"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.913"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.62%  "
